# Updates currentAveragePrice / profit figures across several crafting-leve
# sheets (ALC, ARM, CUL, GSM, LTW, WVR) to reflect refreshed market board
# data pulled by the scheduled runner. Mirrors the authoritative XML diff
# row-for-row, cell-for-cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1086.6666
$ws.Range("I70").Value = 947.5
$ws.Range("J70").Value = 1365
$ws.Range("K70").Value = 2842.5
$ws.Range("L70").Value = 4095
$ws.Range("M70").Value = -2572.5
$ws.Range("N70").Value = -4635

$ws.Range("H73").Value = 1086.6666
$ws.Range("I73").Value = 947.5
$ws.Range("J73").Value = 1365
$ws.Range("K73").Value = 2842.5
$ws.Range("L73").Value = 4095
$ws.Range("M73").Value = -1906.5
$ws.Range("N73").Value = -5967

$ws.Range("H86").Value = 2632.0833
$ws.Range("I86").Value = 2550.6
$ws.Range("J86").Value = 2690.2856
$ws.Range("K86").Value = 2550.6
$ws.Range("L86").Value = 2690.2856
$ws.Range("M86").Value = -1427.6
$ws.Range("N86").Value = -4936.2856

$ws.Range("H88").Value = 5267894
$ws.Range("I88").Value = 8762.5
$ws.Range("J88").Value = 9092717
$ws.Range("K88").Value = 8762.5
$ws.Range("L88").Value = 9092717
$ws.Range("M88").Value = -8356.5
$ws.Range("N88").Value = -9093529

$ws.Range("H89").Value = 2632.0833
$ws.Range("I89").Value = 2550.6
$ws.Range("J89").Value = 2690.2856
$ws.Range("K89").Value = 12753
$ws.Range("L89").Value = 13451.428
$ws.Range("M89").Value = -7137
$ws.Range("N89").Value = -24683.428

$ws.Range("H91").Value = 5267894
$ws.Range("I91").Value = 8762.5
$ws.Range("J91").Value = 9092717
$ws.Range("K91").Value = 8762.5
$ws.Range("L91").Value = 9092717
$ws.Range("M91").Value = -7358.5
$ws.Range("N91").Value = -9095525

$ws.Range("H93").Value = 33067.332
$ws.Range("J93").Value = 33067.332
$ws.Range("L93").Value = 33067.332
$ws.Range("N93").Value = -38059.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9776.93
$ws.Range("I32").Value = 11027.372
$ws.Range("J32").Value = 5936.2856
$ws.Range("K32").Value = 11027.372
$ws.Range("L32").Value = 5936.2856
$ws.Range("M32").Value = -10740.372
$ws.Range("N32").Value = -6510.2856

$ws.Range("H88").Value = 2338.4546
$ws.Range("I88").Value = 1433.3334
$ws.Range("J88").Value = 2677.875
$ws.Range("K88").Value = 1433.3334
$ws.Range("L88").Value = 2677.875
$ws.Range("M88").Value = -1027.3334
$ws.Range("N88").Value = -3489.875

$ws.Range("H91").Value = 2338.4546
$ws.Range("I91").Value = 1433.3334
$ws.Range("J91").Value = 2677.875
$ws.Range("K91").Value = 1433.3334
$ws.Range("L91").Value = 2677.875
$ws.Range("M91").Value = -29.33339999999998
$ws.Range("N91").Value = -5485.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 874.2143
$ws.Range("I117").Value = 475.66666
$ws.Range("K117").Value = 1426.99998
$ws.Range("M117").Value = 2015.00002

$ws.Range("H129").Value = 3540.3872
$ws.Range("I129").Value = 5292.727
$ws.Range("J129").Value = 2576.6
$ws.Range("K129").Value = 15878.181
$ws.Range("L129").Value = 7729.799999999999
$ws.Range("M129").Value = -10878.181
$ws.Range("N129").Value = -17729.8

$ws.Range("H131").Value = 826.4400000000001
$ws.Range("I131").Value = 460.875
$ws.Range("J131").Value = 896.0714
$ws.Range("K131").Value = 1382.625
$ws.Range("L131").Value = 2688.2142
$ws.Range("M131").Value = 3657.375
$ws.Range("N131").Value = -12768.2142

$ws.Range("H134").Value = 3937.3076
$ws.Range("I134").Value = 2673.125
$ws.Range("J134").Value = 5960
$ws.Range("K134").Value = 8019.375
$ws.Range("L134").Value = 17880
$ws.Range("M134").Value = -2949.375
$ws.Range("N134").Value = -28020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 3800
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 3800
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3800
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -4080

$ws.Range("H80").Value = 12348399
$ws.Range("I80").Value = 27780346
$ws.Range("J80").Value = 2842.2666
$ws.Range("K80").Value = 27780346
$ws.Range("L80").Value = 2842.2666
$ws.Range("M80").Value = -27779348
$ws.Range("N80").Value = -4838.2666

$ws.Range("H83").Value = 12348399
$ws.Range("I83").Value = 27780346
$ws.Range("J83").Value = 2842.2666
$ws.Range("K83").Value = 138901730
$ws.Range("L83").Value = 14211.333
$ws.Range("M83").Value = -138896738
$ws.Range("N83").Value = -24195.333

$ws.Range("H132").Value = 6523.7144
$ws.Range("I132").Value = 5601.2
$ws.Range("J132").Value = 7036.222
$ws.Range("K132").Value = 16803.6
$ws.Range("L132").Value = 21108.666
$ws.Range("M132").Value = -14273.6
$ws.Range("N132").Value = -26168.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1013.6
$ws.Range("I22").Value = 660.9
$ws.Range("J22").Value = 1248.7333
$ws.Range("K22").Value = 660.9
$ws.Range("L22").Value = 1248.7333
$ws.Range("M22").Value = -365.9
$ws.Range("N22").Value = -1838.7333

$ws.Range("H27").Value = 1013.6
$ws.Range("I27").Value = 660.9
$ws.Range("J27").Value = 1248.7333
$ws.Range("K27").Value = 660.9
$ws.Range("L27").Value = 1248.7333
$ws.Range("M27").Value = -553.9
$ws.Range("N27").Value = -1462.7333

$ws.Range("H82").Value = 1858.5
$ws.Range("I82").Value = 1570.3334
$ws.Range("J82").Value = 2146.6667
$ws.Range("K82").Value = 1570.3334
$ws.Range("L82").Value = 2146.6667
$ws.Range("M82").Value = -1209.3334
$ws.Range("N82").Value = -2868.6667

$ws.Range("H85").Value = 1858.5
$ws.Range("I85").Value = 1570.3334
$ws.Range("J85").Value = 2146.6667
$ws.Range("K85").Value = 1570.3334
$ws.Range("L85").Value = 2146.6667
$ws.Range("M85").Value = -322.3334
$ws.Range("N85").Value = -4642.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16286.143
$ws.Range("J62").Value = 16286.143
$ws.Range("L62").Value = 16286.143
$ws.Range("N62").Value = -17534.143

$ws.Range("H65").Value = 16286.143
$ws.Range("J65").Value = 16286.143
$ws.Range("L65").Value = 81430.715
$ws.Range("N65").Value = -87670.715

$ws.Range("H81").Value = 1776.421
$ws.Range("I81").Value = 1065.7142
$ws.Range("J81").Value = 2191
$ws.Range("K81").Value = 2131.4284
$ws.Range("L81").Value = 4382
$ws.Range("M81").Value = -1070.4284
$ws.Range("N81").Value = -6504

$ws.Range("H84").Value = 1776.421
$ws.Range("I84").Value = 1065.7142
$ws.Range("J84").Value = 2191
$ws.Range("K84").Value = 10657.142
$ws.Range("L84").Value = 21910
$ws.Range("M84").Value = -5353.142
$ws.Range("N84").Value = -32518

$ws.Range("H92").Value = 13000
$ws.Range("J92").Value = 13000
$ws.Range("L92").Value = 13000
$ws.Range("N92").Value = -17992

$ws.Range("H96").Value = 4227.846
$ws.Range("I96").Value = 2261
$ws.Range("J96").Value = 5913.7144
$ws.Range("K96").Value = 2261
$ws.Range("L96").Value = 5913.7144
$ws.Range("M96").Value = -888
$ws.Range("N96").Value = -8659.714400000001
